$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue $ws "D2" "27.194.39"
Set-TextValue $ws "E2" "  -0.02%  "

Set-TextValue $ws "D3" "1.903.10"
Set-TextValue $ws "E3" "  -0.06%  "

Set-TextValue $ws "D4" "1.003"
Set-TextValue $ws "E4" "  +0.25%  "

Set-TextValue $ws "D5" "307.48"
Set-TextValue $ws "E5" "  +0.28%  "

Set-TextValue $ws "E6" "  +0.25%  "

Set-TextValue $ws "D7" "0.5248"
Set-TextValue $ws "E7" "  -0.30%  "

Set-TextValue $ws "D8" "0.3814"
Set-TextValue $ws "E8" "  +0.99%  "

Set-TextValue $ws "D9" "0.07300"
Set-TextValue $ws "E9" "  +0.63%  "

Set-TextValue $ws "D10" "21.58"
Set-TextValue $ws "E10" "  +2.10%  "

Set-TextValue $ws "D11" "0.9051"
Set-TextValue $ws "E11" "  +0.58%  "

Set-TextValue $ws "D12" "0.08132"
Set-TextValue $ws "E12" "  -3.46%  "

Set-TextValue $ws "D13" "95.60"
Set-TextValue $ws "E13" "  +0.92%  "

Set-TextValue $ws "D14" "5.353"
Set-TextValue $ws "E14" "  +1.55%  "

Set-TextValue $ws "D15" "1.811.93"
Set-TextValue $ws "E15" "  -4.86%  "

Set-TextValue $ws "E16" "  +0.25%  "

Set-TextValue $ws "D17" "0.000008661"
Set-TextValue $ws "E17" "  +0.48%  "

Set-TextValue $ws "D18" "14.71"
Set-TextValue $ws "E18" "  +0.98%  "

Set-TextValue $ws "D19" "1.001"
Set-TextValue $ws "E19" "  +0.17%  "

Set-TextValue $ws "D20" "27.228.23"
Set-TextValue $ws "E20" "  +0.00%  "

Set-TextValue $ws "D21" "5.110"
Set-TextValue $ws "E21" "  +0.98%  "

Set-TextValue $ws "E22" "  +2.04%  "

Set-TextValue $ws "D23" "6.461"
Set-TextValue $ws "E23" "  +0.37%  "

Set-TextValue $ws "D24" "2.328"
Set-TextValue $ws "E24" "  +2.31%  "

Set-TextValue $ws "D25" "149.30"
Set-TextValue $ws "E25" "  +1.64%  "

Set-TextValue $ws "D26" "18.24"
Set-TextValue $ws "E26" "  +0.45%  "

Set-TextValue $ws "D27" "1.735"
Set-TextValue $ws "E27" "  -0.88%  "

Set-TextValue $ws "D28" "116.13"
Set-TextValue $ws "E28" "  +1.09%  "

Set-TextValue $ws "D29" "4.832"
Set-TextValue $ws "E29" "  +0.39%  "

Set-TextValue $ws "D30" "4.888"
Set-TextValue $ws "E30" "  -0.69%  "

Set-TextValue $ws "D31" "0.09237"
Set-TextValue $ws "E31" "  -0.62%  "

Set-TextValue $ws "D32" "0.05073"
Set-TextValue $ws "E32" "  +0.18%  "

Set-TextValue $ws "D33" "0.7953"
Set-TextValue $ws "E33" "  -1.50%  "

Set-TextValue $ws "D34" "1.228"
Set-TextValue $ws "E34" "  -0.73%  "

Set-TextValue $ws "E35" "  +0.86%  "

Set-TextValue $ws "D36" "3.367"
Set-TextValue $ws "E36" "  -0.07%  "

Set-TextValue $ws "D37" "2.662"
Set-TextValue $ws "E37" "  +1.45%  "

Set-TextValue $ws "D38" "0.5718"
Set-TextValue $ws "E38" "  -0.13%  "

Set-TextValue $ws "D39" "0.01991"
Set-TextValue $ws "E39" "  +0.18%  "

Set-TextValue $ws "D40" "1.084"
Set-TextValue $ws "E40" "  +0.99%  "

Set-TextValue $ws "D41" "9.029"
Set-TextValue $ws "E41" "  +0.71%  "

Set-TextValue $ws "D42" "6.589"
Set-TextValue $ws "E42" "  -0.81%  "

Set-TextValue $ws "D43" "116.23"
Set-TextValue $ws "E43" "  -1.17%  "

Set-TextValue $ws "D44" "0.1514"
Set-TextValue $ws "E44" "  -0.13%  "

Set-TextValue $ws "D45" "0.4882"
Set-TextValue $ws "E45" "  +0.80%  "

# Row 46/47: EnergySwap and PaxDollar swap places (Coin name + Link), with refreshed Price/Volume
Set-TextValue $ws "B46" "EnergySwap"
Set-TextValue $ws "C46" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws "D46" "10.19"
Set-TextValue $ws "E46" "  +0.20%  "

Set-TextValue $ws "B47" "PaxDollar"
Set-TextValue $ws "C47" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws "D47" "1.003"
Set-TextValue $ws "E47" "  +0.31%  "


Set-TextValue $ws "D48" "1.633"
Set-TextValue $ws "E48" "  +1.01%  "

Set-TextValue $ws "D49" "38.54"
Set-TextValue $ws "E49" "  +2.97%  "

Set-TextValue $ws "D50" "64.00"
Set-TextValue $ws "E50" "  +0.51%  "

Set-TextValue $ws "D51" "0.05959"
Set-TextValue $ws "E51" "  +0.38%  "
